$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty row 12, shifting rows 13:24 up to 12:23
$ws.Rows.Item(12).Delete()

# Select entire row 12 (matches the recorded selection after the delete)
$ws.Range("A12:XFD12").Select()
